# Updated cryptos list on Wed Aug 23 23:57:20 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple Price (D) / Volume(1h) (E) updates for rows 2-44 ---

$ws.Range("D2").Value = "'26.463.83"
$ws.Range("E2").Value = "'  +1.57%  "

$ws.Range("D3").Value = "'1.680.72"
$ws.Range("E3").Value = "'  +2.79%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "'  +0.15%  "

$ws.Range("D5").Value = "'216.58"
$ws.Range("E5").Value = "'  +2.66%  "

$ws.Range("D6").Value = "'0.5328"
$ws.Range("E6").Value = "'  +1.84%  "

$ws.Range("E7").Value = "'  +0.08%  "

$ws.Range("D8").Value = "'0.2688"
$ws.Range("E8").Value = "'  +3.79%  "

$ws.Range("D9").Value = "'0.06404"
$ws.Range("E9").Value = "'  +2.00%  "

$ws.Range("E10").Value = "'  +5.57%  "

$ws.Range("D11").Value = "'0.07805"
$ws.Range("E11").Value = "'  +2.93%  "

$ws.Range("D12").Value = "'1.683.18"
$ws.Range("E12").Value = "'  +3.15%  "

$ws.Range("D13").Value = "'4.495"
$ws.Range("E13").Value = "'  +1.55%  "

$ws.Range("D14").Value = "'0.5579"
$ws.Range("E14").Value = "'  +1.56%  "

$ws.Range("D15").Value = "'0.0₅8325"
$ws.Range("E15").Value = "'  +4.02%  "

$ws.Range("D16").Value = "'65.67"
$ws.Range("E16").Value = "'  +1.39%  "

$ws.Range("D17").Value = "'26.520.35"
$ws.Range("E17").Value = "'  +1.90%  "

$ws.Range("E18").Value = "'  +0.00%  "

$ws.Range("D19").Value = "'4.763"
$ws.Range("E19").Value = "'  +1.96%  "

$ws.Range("D20").Value = "'194.64"
$ws.Range("E20").Value = "'  +4.93%  "

$ws.Range("E21").Value = "'  +2.13%  "

$ws.Range("D22").Value = "'6.358"
$ws.Range("E22").Value = "'  +3.90%  "

$ws.Range("E23").Value = "'  +0.16%  "

$ws.Range("D24").Value = "'143.20"
$ws.Range("E24").Value = "'  -1.61%  "

$ws.Range("D25").Value = "'0.1278"
$ws.Range("E25").Value = "'  +5.46%  "

$ws.Range("D26").Value = "'7.440"
$ws.Range("E26").Value = "'  +0.46%  "

$ws.Range("D27").Value = "'16.35"
$ws.Range("E27").Value = "'  +4.70%  "

$ws.Range("D28").Value = "'1.427"
$ws.Range("E28").Value = "'  +3.67%  "

$ws.Range("D29").Value = "'0.06227"
$ws.Range("E29").Value = "'  +4.89%  "

$ws.Range("D30").Value = "'1.273"
$ws.Range("E30").Value = "'  +2.44%  "

$ws.Range("D31").Value = "'3.606"
$ws.Range("E31").Value = "'  +5.20%  "

$ws.Range("E32").Value = "'  +1.64%  "

$ws.Range("D33").Value = "'1.689"
$ws.Range("E33").Value = "'  +3.51%  "

$ws.Range("E34").Value = "'  +2.85%  "

$ws.Range("D35").Value = "'2.428"
$ws.Range("E35").Value = "'  +1.65%  "

$ws.Range("D36").Value = "'2.789"
$ws.Range("E36").Value = "'  +1.47%  "

$ws.Range("D37").Value = "'0.5749"
$ws.Range("E37").Value = "'  -0.63%  "

$ws.Range("D38").Value = "'0.01637"
$ws.Range("E38").Value = "'  +2.18%  "

$ws.Range("D39").Value = "'6.048"
$ws.Range("E39").Value = "'  +6.57%  "

$ws.Range("D40").Value = "'1.074.79"
$ws.Range("E40").Value = "'  +3.32%  "

$ws.Range("D41").Value = "'0.8576"
$ws.Range("E41").Value = "'  +1.16%  "

$ws.Range("D42").Value = "'1.001"

$ws.Range("E43").Value = "'  -0.19%  "

$ws.Range("D44").Value = "'1.827.47"
$ws.Range("E44").Value = "'  +2.49%  "

# --- Rows 45-51: BabyDogeCoin inserted, pushing Aave..Mantle down one row,
#     RenderToken (old row 51) drops off the bottom of the list ---

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "'0.0₈112"
$ws.Range("E45").Value = "'  +4.03%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'57.22"
$ws.Range("E46").Value = "'  +4.17%  "

$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "'  +0.71%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.108"
$ws.Range("E48").Value = "'  +1.13%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05205"
$ws.Range("E49").Value = "'  +0.91%  "

$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "'6.030"
$ws.Range("E50").Value = "'  +3.00%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.4239"
$ws.Range("E51").Value = "'  +0.36%  "

